$d = $word.ActiveDocument

# 1. Update the "Directed a team of developers" bullet
$d.Content.Find.Execute(
    "Directed a team of developers: trained, delegated, instituted standards, and reviewed code",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Directed a team of developers in writing and testing applications for Hotel Dieu Hospital and the Ministry of Health and Long-Term Care",
    2)

# 2. Tweak the "Devised a web UI" bullet wording
$d.Content.Find.Execute(
    "Devised a web UI that simulated the Windows desktop as part of an asset management application, and architected the backing database",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Devised a web UI simulating the Windows desktop as part of an asset management application, and architected the backing database",
    2)

# 3. Remove the now-redundant "Wrote and tested applications..." bullet entirely
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Wrote and tested applications for Hotel Dieu Hospital and the Ministry of Health and Long-Term Care*") {
        $p.Range.Delete()
        break
    }
}
